$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 9 data documenting the "No Account Lockout Mechanism" security issue.
# Cells are written in B9, C9, E9, D9 order so that the workbook's shared-string
# table gains the four new entries in the same order as the target file
# (index 29 = B9 text, 30 = C9 text, 31 = E9 text, 32 = D9 text).
$ws.Range("B9").Value = "No Account Lockout Mechanism "
$ws.Range("C9").Value = "Main.java, Login.Java"
$ws.Range("E9").Value = "A lockout mechanism was implemented where after 5 failed login attempts, the user is blocked from making another attemp for 5 minutes. " + [char]10 + " The validateLogin function was modified to include a failed login attempt counter.  "
$ws.Range("D9").Value = "Brute Force Attacks"

# Reflect the last active/selected cell recorded in the saved workbook.
[void]$ws.Range("D8").Select()
